$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix letter case inconsistencies in header labels (SSDM-12286)
$ws.Range("H4").Value = "Vocabulary code"
$ws.Range("H12").Value = "Vocabulary code"
$ws.Range("E2").Value = "Generated code prefix"
$ws.Range("E10").Value = "Generated code prefix"

$ws.Range("E10").Select()
